# Quarterly indexing esoteric bug-fix operation
#
# Each data row (2..16) on the sheet held the QoQ error series starting at
# column B. The fix shifts every existing value one column to the right
# (B->C, C->D, ... J->K, dropping whatever used to be in K) and writes a
# freshly computed value into the now-vacated column B. Column A (the
# row label) and row 1 (the header) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to place into column B for each data row.
$newB = @{
    2  = -0.2177157015159319
    3  = -0.1395947820665385
    4  = -0.3119065001142551
    5  = 0.7021231295320197
    6  = 1.514070997382048
    7  = 0.2163102553365951
    8  = 0.3684555432821496
    9  = 0.661541622456546
    10 = -0.07992401592518952
    11 = 0.1551026493581833
    12 = -0.08373363042288225
    13 = 0.1925427069667326
    14 = -0.4379379024501944
    15 = 0.2324016585002178
    16 = -0.09587373626955231
}

# Last populated data column (before the edit) for each row; K = 11.
$lastCol = @{
    2  = 11
    3  = 11
    4  = 11
    5  = 11
    6  = 11
    7  = 10
    8  = 9
    9  = 8
    10 = 7
    11 = 6
    12 = 5
    13 = 4
    14 = 3
    15 = 2
    16 = 1
}

for ($r = 2; $r -le 16; $r++) {
    $last = $lastCol[$r]

    # Shift existing values one column to the right, starting from the
    # rightmost column so we never clobber a value before reading it.
    # If the row is already full out to column K (11), the old K value
    # simply falls off the end (matches the source diff).
    $destStart = [Math]::Min($last + 1, 11)
    for ($c = $destStart; $c -ge 3; $c--) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($r, $c - 1).Value2
    }

    # Write the new value into the vacated column B.
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
